# Applies the NIT-9006377433 data refresh:
#  - Row 16-17: replace XAVIER ELIAS TORRES MIRANDA (1810) rows with
#    GORTRUDE MARY MCLEAN CARDILES (1906/1905), Valor Mora 3,000,000
#  - Rows 18-35: keep XAVIER ELIAS TORRES MIRANDA but shift/renumber the
#    "Periodo Mora" sequence and update Salario Basico to 781,242
#    (except the first of that block, which keeps a 20,833 Valor Mora)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(16, 2).Value = "CC"
$ws.Cells.Item(16, 3).Value = "1082925993"
$ws.Cells.Item(16, 4).Value = "GORTRUDE MARY MCLEAN CARDILES"
$ws.Cells.Item(16, 5).Value = "1906"
$ws.Cells.Item(16, 6).Value = 31249
$ws.Cells.Item(16, 7).Value = 3000000
$ws.Cells.Item(17, 2).Value = "CC"
$ws.Cells.Item(17, 3).Value = "1082925993"
$ws.Cells.Item(17, 4).Value = "GORTRUDE MARY MCLEAN CARDILES"
$ws.Cells.Item(17, 5).Value = "1905"
$ws.Cells.Item(17, 6).Value = 31249
$ws.Cells.Item(17, 7).Value = 3000000
$ws.Cells.Item(18, 2).Value = "CC"
$ws.Cells.Item(18, 3).Value = "73009373"
$ws.Cells.Item(18, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(18, 5).Value = "2003"
$ws.Cells.Item(18, 6).Value = 20833
$ws.Cells.Item(18, 7).Value = 781242
$ws.Cells.Item(19, 2).Value = "CC"
$ws.Cells.Item(19, 3).Value = "73009373"
$ws.Cells.Item(19, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(19, 5).Value = "2002"
$ws.Cells.Item(19, 6).Value = 31249
$ws.Cells.Item(19, 7).Value = 781242
$ws.Cells.Item(20, 2).Value = "CC"
$ws.Cells.Item(20, 3).Value = "73009373"
$ws.Cells.Item(20, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(20, 5).Value = "2001"
$ws.Cells.Item(20, 6).Value = 31249
$ws.Cells.Item(20, 7).Value = 781242
$ws.Cells.Item(21, 2).Value = "CC"
$ws.Cells.Item(21, 3).Value = "73009373"
$ws.Cells.Item(21, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(21, 5).Value = "1912"
$ws.Cells.Item(21, 6).Value = 31249
$ws.Cells.Item(21, 7).Value = 781242
$ws.Cells.Item(22, 2).Value = "CC"
$ws.Cells.Item(22, 3).Value = "73009373"
$ws.Cells.Item(22, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(22, 5).Value = "1911"
$ws.Cells.Item(22, 6).Value = 31249
$ws.Cells.Item(22, 7).Value = 781242
$ws.Cells.Item(23, 2).Value = "CC"
$ws.Cells.Item(23, 3).Value = "73009373"
$ws.Cells.Item(23, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(23, 5).Value = "1910"
$ws.Cells.Item(23, 6).Value = 31249
$ws.Cells.Item(23, 7).Value = 781242
$ws.Cells.Item(24, 2).Value = "CC"
$ws.Cells.Item(24, 3).Value = "73009373"
$ws.Cells.Item(24, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(24, 5).Value = "1909"
$ws.Cells.Item(24, 6).Value = 31249
$ws.Cells.Item(24, 7).Value = 781242
$ws.Cells.Item(25, 2).Value = "CC"
$ws.Cells.Item(25, 3).Value = "73009373"
$ws.Cells.Item(25, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(25, 5).Value = "1908"
$ws.Cells.Item(25, 6).Value = 31249
$ws.Cells.Item(25, 7).Value = 781242
$ws.Cells.Item(26, 2).Value = "CC"
$ws.Cells.Item(26, 3).Value = "73009373"
$ws.Cells.Item(26, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(26, 5).Value = "1907"
$ws.Cells.Item(26, 6).Value = 31249
$ws.Cells.Item(26, 7).Value = 781242
$ws.Cells.Item(27, 2).Value = "CC"
$ws.Cells.Item(27, 3).Value = "73009373"
$ws.Cells.Item(27, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(27, 5).Value = "1906"
$ws.Cells.Item(27, 6).Value = 31249
$ws.Cells.Item(27, 7).Value = 781242
$ws.Cells.Item(28, 2).Value = "CC"
$ws.Cells.Item(28, 3).Value = "73009373"
$ws.Cells.Item(28, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(28, 5).Value = "1905"
$ws.Cells.Item(28, 6).Value = 31249
$ws.Cells.Item(28, 7).Value = 781242
$ws.Cells.Item(29, 2).Value = "CC"
$ws.Cells.Item(29, 3).Value = "73009373"
$ws.Cells.Item(29, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(29, 5).Value = "1904"
$ws.Cells.Item(29, 6).Value = 31249
$ws.Cells.Item(29, 7).Value = 781242
$ws.Cells.Item(30, 2).Value = "CC"
$ws.Cells.Item(30, 3).Value = "73009373"
$ws.Cells.Item(30, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(30, 5).Value = "1903"
$ws.Cells.Item(30, 6).Value = 31249
$ws.Cells.Item(30, 7).Value = 781242
$ws.Cells.Item(31, 2).Value = "CC"
$ws.Cells.Item(31, 3).Value = "73009373"
$ws.Cells.Item(31, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(31, 5).Value = "1902"
$ws.Cells.Item(31, 6).Value = 31249
$ws.Cells.Item(31, 7).Value = 781242
$ws.Cells.Item(32, 2).Value = "CC"
$ws.Cells.Item(32, 3).Value = "73009373"
$ws.Cells.Item(32, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(32, 5).Value = "1901"
$ws.Cells.Item(32, 6).Value = 31249
$ws.Cells.Item(32, 7).Value = 781242
$ws.Cells.Item(33, 2).Value = "CC"
$ws.Cells.Item(33, 3).Value = "73009373"
$ws.Cells.Item(33, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(33, 5).Value = "1812"
$ws.Cells.Item(33, 6).Value = 31249
$ws.Cells.Item(33, 7).Value = 781242
$ws.Cells.Item(34, 2).Value = "CC"
$ws.Cells.Item(34, 3).Value = "73009373"
$ws.Cells.Item(34, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(34, 5).Value = "1811"
$ws.Cells.Item(34, 6).Value = 31249
$ws.Cells.Item(34, 7).Value = 781242
$ws.Cells.Item(35, 2).Value = "CC"
$ws.Cells.Item(35, 3).Value = "73009373"
$ws.Cells.Item(35, 4).Value = "XAVIER ELIAS TORRES MIRANDA"
$ws.Cells.Item(35, 5).Value = "1810"
$ws.Cells.Item(35, 6).Value = 31249
$ws.Cells.Item(35, 7).Value = 781242

